$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.254102945327759
$ws.Range("B1").Value = 2.212313175201416
$ws.Range("C1").Value = 4.334763526916504
$ws.Range("D1").Value = 3.036949396133423
$ws.Range("E1").Value = 1.040203452110291
